# Data provider for Car search added
# Adds DropOff/PickUp date values (as text, "mm/dd/yyyy" strings) for the
# three sample rows, centers all the data, and widens a couple of columns
# to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- PickUp / DropOff dates, entered as centered text values -------------
# (NumberFormat is set to "@" BEFORE the value so Excel stores them as
# shared-string text instead of coercing the mm/dd/yyyy text into a date
# serial number.)

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "07/20/2020"
$ws.Range("B2").HorizontalAlignment = $xlCenter

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "07/25/2020"
$ws.Range("B3").HorizontalAlignment = $xlCenter

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "07/22/2020"
$ws.Range("B4").HorizontalAlignment = $xlCenter

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "08/30/2020"
$ws.Range("D2").HorizontalAlignment = $xlCenter

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "08/30/2020"
$ws.Range("D3").HorizontalAlignment = $xlCenter

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "07/30/2020"
$ws.Range("D4").HorizontalAlignment = $xlCenter

# --- Center the rest of the table (headers + location/time columns) ------
$ws.Range("A1:E4").HorizontalAlignment = $xlCenter

# --- Column widths, to fit the wider date text ----------------------------
$ws.Range("B1").ColumnWidth = 14.25
$ws.Range("C1").ColumnWidth = 15.42
$ws.Range("F1").ColumnWidth = 11.25

# --- Active cell/selection, as left by the editing session ---------------
$ws.Range("C8").Select() | Out-Null
